# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-10-21 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-10-22 Wednesday", 2)

# Update the multiplication problems in the table. The table contains
# several duplicate "op=" expressions, so address each cell by its
# (row, column) position rather than relying on text search alone.
$t = $d.Tables.Item(1)

function Set-Cell($table, $row, $col, $oldText, $newText) {
    $cell = $table.Cell($row, $col)
    $current = $cell.Range.Text
    if ($current -notlike "$oldText*") {
        Write-Host "WARNING: cell ($row,$col) expected '$oldText' but found '$current'"
    }
    $cell.Range.Text = $newText
}

Set-Cell $t 1 1 "35×25=" "97×33="
Set-Cell $t 1 2 "61×55=" "34×95="
Set-Cell $t 1 3 "35×88=" "67×13="
Set-Cell $t 1 4 "45×35=" "66×31="
Set-Cell $t 1 5 "54×35=" "72×80="

Set-Cell $t 5 1 "66×37=" "25×37="
Set-Cell $t 5 2 "58×26=" "93×60="
Set-Cell $t 5 3 "44×29=" "49×26="
Set-Cell $t 5 4 "80×93=" "47×28="
Set-Cell $t 5 5 "28×90=" "84×32="

Set-Cell $t 10 1 "76×20=" "39×16="
Set-Cell $t 10 2 "16×74=" "85×78="
Set-Cell $t 10 3 "19×43=" "63×51="
Set-Cell $t 10 4 "80×93=" "21×36="
Set-Cell $t 10 5 "30×69=" "21×98="

Set-Cell $t 15 1 "37×50=" "98×54="
Set-Cell $t 15 2 "13×27=" "50×36="
Set-Cell $t 15 3 "43×50=" "47×12="
Set-Cell $t 15 4 "36×42=" "82×32="
Set-Cell $t 15 5 "30×21=" "30×87="

Set-Cell $t 20 1 "69×84=" "44×49="
Set-Cell $t 20 2 "11×85=" "58×20="
Set-Cell $t 20 3 "81×35=" "61×20="
Set-Cell $t 20 4 "56×35=" "26×82="
Set-Cell $t 20 5 "28×33=" "66×24="

Write-Host "Done"
